# Updated symbol list on Sun Feb  5 18:46:46 UTC 2023 with GitHub Actions
# Applies the refreshed cryptocurrency price / volume figures (and the
# GateToken <-> FTXToken row swap) to the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free approach: set the Price (D) and Volume(1h) (E) columns as
# plain text so Excel doesn't silently reinterpret numeric- or
# percent-looking strings as actual numbers (which would lose the exact
# textual formatting, e.g. trailing zeros, used in the source data).
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

# Row 2 - BNB
$ws.Range("D2").Value = "328.39"
$ws.Range("E2").Value = "-0.82%"

# Row 3 - OKB
$ws.Range("D3").Value = "43.93"
$ws.Range("E3").Value = "5.58%"

# Row 4 - HuobiToken
$ws.Range("D4").Value = "5.412"
$ws.Range("E4").Value = "-4.81%"

# Row 5 - Cronos
$ws.Range("D5").Value = "0.08099"
$ws.Range("E5").Value = "-2.97%"

# Row 6 - KuCoinToken
$ws.Range("D6").Value = "8.688"
$ws.Range("E6").Value = "-1.36%"

# Row 7 & Row 8 swap places: GateToken <-> FTXToken
# Row 7 becomes FTXToken
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "1.905"
$ws.Range("E7").Value = "-4.74%"

# Row 8 becomes GateToken
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "4.304"
$ws.Range("E8").Value = "-3.66%"

# Row 9 - BTSEToken
$ws.Range("D9").Value = "2.751"
$ws.Range("E9").Value = "-5.00%"

# Row 10 - MXToken
$ws.Range("D10").Value = "0.9439"
$ws.Range("E10").Value = "1.94%"

# Row 11 - LiechtensteinCryptoassetsExchange
$ws.Range("D11").Value = "0.1180"
$ws.Range("E11").Value = "-8.24%"

# Row 12 - WazirX
$ws.Range("D12").Value = "0.1894"
$ws.Range("E12").Value = "-4.39%"

# Row 13 - MandalaExchangeToken
$ws.Range("D13").Value = "0.09566"
$ws.Range("E13").Value = "0.73%"

# Row 14 - BitrueCoin
$ws.Range("D14").Value = "0.04190"
$ws.Range("E14").Value = "8.94%"

# Row 15 - BitMartToken
$ws.Range("D15").Value = "0.1070"
$ws.Range("E15").Value = "0.92%"

# Row 16 - BitForexToken
$ws.Range("D16").Value = "0.001284"
$ws.Range("E16").Value = "-1.62%"

# Row 17 - TigerCash
$ws.Range("D17").Value = "0.005948"
$ws.Range("E17").Value = "-2.65%"

# Row 18 - LEO
$ws.Range("D18").Value = "3.554"
$ws.Range("E18").Value = "3.46%"

# Row 20 - MCDex
$ws.Range("D20").Value = "8.517"
$ws.Range("E20").Value = "-1.70%"

# Row 21 - ProBitToken
$ws.Range("D21").Value = "0.1361"
$ws.Range("E21").Value = "-0.15%"

# Row 22 - ZBToken
$ws.Range("D22").Value = "0.2608"
$ws.Range("E22").Value = "4.93%"

# Row 23 - CoinExToken
$ws.Range("D23").Value = "0.04390"
$ws.Range("E23").Value = "-0.64%"

# Row 24 - BitKan (price unchanged, only volume)
$ws.Range("E24").Value = "-2.77%"

# Row 25 - HotbitToken
$ws.Range("D25").Value = "0.004311"
$ws.Range("E25").Value = "-1.79%"

# Row 26 - NitroEx (price unchanged, only volume)
$ws.Range("E26").Value = "1.42%"

# Row 27 - UpBots
$ws.Range("D27").Value = "0.0004021"
$ws.Range("E27").Value = "32.06%"

# Row 39 - One
$ws.Range("D39").Value = "0.02690"
$ws.Range("E39").Value = "-4.61%"

# Row 40 - IDEX
$ws.Range("D40").Value = "0.05498"
$ws.Range("E40").Value = "-0.65%"

# Row 41 - KickToken
$ws.Range("D41").Value = "0.007810"
$ws.Range("E41").Value = "-1.73%"

# Row 42 - Dexo
$ws.Range("D42").Value = "0.009785"
$ws.Range("E42").Value = "5.22%"

# Row 43 - BKEXToken
$ws.Range("D43").Value = "0.1394"
$ws.Range("E43").Value = "-2.69%"

# Row 44 - CEJI
$ws.Range("D44").Value = "0.002127"
$ws.Range("E44").Value = "0.72%"

# Row 45 - LocalTraders
$ws.Range("D45").Value = "0.009632"
$ws.Range("E45").Value = "-17.90%"

# Row 46 - CoinLion
$ws.Range("D46").Value = "0.00007118"
$ws.Range("E46").Value = "2.68%"

# Row 47 - Kangarootoken
$ws.Range("D47").Value = "0.00000000756"
$ws.Range("E47").Value = "0.72%"

# Row 48 - BOLO
$ws.Range("D48").Value = "0.003477"
$ws.Range("E48").Value = "0.80%"

# Row 49 - CoinbaseStockToken
$ws.Range("D49").Value = "0.002288"
$ws.Range("E49").Value = "0.38%"

# Row 50 - CryptobidCoin
$ws.Range("D50").Value = "0.00002117"
$ws.Range("E50").Value = "0.72%"

# Row 51 - SpecialPowerGold
$ws.Range("D51").Value = "0.0002016"
$ws.Range("E51").Value = "0.72%"

Write-Output "Applied cryptos.xlsx price/volume refresh"
